$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in NAT portforwarding task description (row 10)
$ws.Range("B10").Value = "NAT portforwarding en kleine aanpassingen in netwerkscripts toegevoegd"

# Fill in new task rows for the test day entries
$ws.Range("B12").Value = "Fysieke testdag op school, focus op de totale test voor demo morgen en testen van laatste aanpassingen + test met ACLs actief. "
$ws.Range("C12").Value = "Thomas en Jelle"

$ws.Range("B13").Value = "CA uitbreiding poging "
$ws.Range("C13").Value = "Jelle"

# Update the active selection to match the author's final cursor position
$ws.Range("B10").Select()
